$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10: 220-224 / Suzano Imperatriz
$ws.Range("A10").Value2 = "220-224"
$ws.Range("B10").Value2 = "Suzano Imperatriz"
$ws.Range("C10").Value2 = $ws.Range("C2").Value2
$ws.Range("C10").Style = $ws.Range("C2").Style
$ws.Range("D10").Value2 = 0

# Row 11: 60 / Brascabos
$ws.Range("A11").Value2 = 60
$ws.Range("B11").Value2 = "Brascabos"
$ws.Range("C11").Value2 = $ws.Range("C2").Value2
$ws.Range("C11").Style = $ws.Range("C2").Style
$ws.Range("D11").Value2 = 0

# Update the selection to match the author's final cursor position
$ws.Range("D12").Select()
